$wb = $excel.ActiveWorkbook

# --- Style changes (both sheets share the same title/header layout) ---
# The title (row 1) previously used a dedicated bold 14pt font; it now
# shares the same bold font as the header row, which itself is recolored
# to white so it reads clearly against its dark-blue fill.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count

    # Title row: drop the explicit 14pt size (keeps existing bold), make white.
    $titleCell = $ws.Cells.Item(1, 1)
    $titleCell.Font.Size = 11
    $titleCell.Font.Color = 16777215

    # Header row: already bold, just make white (sits on the dark-blue fill).
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Color = 16777215
}

# --- Data changes on the "Training Dashboard" sheet ---
$ws1 = $wb.Worksheets.Item("Training Dashboard")

# PERIOD TO EXPIRE
$ws1.Range("H3").Value = -50

# LAST UPDATE - enter via formula then flatten to a literal value so the
# date-like text isn't auto-converted into a real date serial number.
$ws1.Range("I3").Formula = '="16-Sep-2025"'
$ws1.Range("I3").Copy() | Out-Null
$ws1.Range("I3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
